$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C4's dialogue text to wrap the parenthetical remark in a <size=24> rich-text tag,
# matching Unity's built-in rich text tag support.
$ws.Range("C4").Value = "Ribbit <i>Ribbit!</i> <size=24>(Yeah <color=green>frog-face!</color> Wrong part of town!)</size>"

# Move the active cell selection to C4 (reflects the last-edited cell).
$ws.Range("C4").Select()
